$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B11 currently holds the text "R40". The cell should now hold the text
# "1" (kept as text, not auto-converted to a number), with its existing
# number format / style left untouched.
$ws.Range("Z1").NumberFormat = "@"
$ws.Range("Z1").Value = "1"
$ws.Range("Z1").Copy()
$ws.Range("B11").PasteSpecial(-4163)   # xlPasteValues - keeps B11's own formatting
$ws.Columns.Item(26).Delete()          # remove helper column Z, restores original dimension
